$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.643.93"
$ws.Range("E2").Value = "  +2.68%  "
$ws.Range("D3").Value = "1.686.47"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5343"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2679"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06434"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07798"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.17%  "
$ws.Range("D12").Value = "1.699.50"
$ws.Range("E12").Value = "  +4.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.501"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5608"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.51%  "
$ws.Range("D15").Value = "0.0₅8454"
$ws.Range("E15").Value = "  +6.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "26.681.17"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.002"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.793"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "195.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.85%  "
$ws.Range("E22").Value = "  +5.15%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1282"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.481"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.434"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06154"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.279"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.608"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.470"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.701"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.012"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.42%  "
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5744"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("E38").Value = "  +3.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.035"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.83%  "
$ws.Range("D40").Value = "1.070.51"
$ws.Range("E40").Value = "  +6.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8635"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.61%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("D44").Value = "1.837.23"
$ws.Range("E44").Value = "  +3.07%  "
$ws.Range("D45").Value = "0.0₈109"
$ws.Range("E45").Value = "  +2.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.188"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.05%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05217"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("E50").Value = "  +5.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4244"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.53%  "
